# Weekly driver report update for 2025-05-05
# Updates the "Good Drivers" table (rows 12-33) on the active sheet with
# refreshed per-driver roaming statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=12; A="Intel(R) Wi-Fi 6 AX201 160MHz - 22.250.0.4"; B=1869842; C=3420;  D=4386; E=2436; F=1877648; H="22.250.0.4";  I=99.59999999999999; J="2023-07-25" },
    @{ Row=13; A="Intel(R) Wi-Fi 6 AX201 160MHz - 22.40.0.7";  B=8170878; C=15867; D=13188; E=18436; F=8199933; H="22.40.0.7";  I=99.59999999999999; J="2021-09-18" },
    @{ Row=14; A="Intel(R) Wi-Fi 6 AX201 160MHz - 22.190.0.4"; B=1611822; C=4793;  D=1136;  E=2996; F=1617751; H="22.190.0.4"; I=99.59999999999999; J="2022-11-22" },
    @{ Row=15; A="Intel(R) Wi-Fi 6 AX201 160MHz - 22.150.1.1"; B=298304;  C=925;   D=303;   E=420;  F=299532;  H="22.150.1.1"; I=99.59999999999999; J="2022-06-20" },
    @{ Row=16; A="Intel(R) Wi-Fi 6 AX201 160MHz - 23.10.0.8";  B=381616;  C=1106;  D=582;   E=601;  F=383304;  H="23.10.0.8";  I=99.59999999999999; J="2023-10-30" },
    @{ Row=17; A="Intel(R) Wi-Fi 6 AX201 160MHz - 22.90.0.5";  B=335610;  C=495;   D=530;   E=414;  F=336635;  H="22.90.0.5";  I=99.7;              J="2021-09-26" },
    @{ Row=18; A="Intel(R) Wi-Fi 6 AX201 160MHz - 22.220.0.4"; B=750778;  C=1655;  D=507;   E=1290; F=752940;  H="22.220.0.4"; I=99.7;              J="2023-03-28" },
    @{ Row=19; A="Intel(R) Wi-Fi 6 AX201 160MHz - 22.20.0.6";  B=1021828; C=2379;  D=1006;  E=1021; F=1025213; H="22.20.0.6";  I=99.7;              J="2020-11-29" },
    @{ Row=20; A="Intel(R) Wi-Fi 6 AX201 160MHz - 22.60.0.6";  B=100380;  C=264;   D=30;    E=267;  F=100674;  H="22.60.0.6";  I=99.7;              J="2021-05-26" },
    @{ Row=21; A="Intel(R) Wi-Fi 6 AX201 160MHz - 23.80.1.3";  B=195464;  C=430;   D=70;    E=401;  F=195964;  H="23.80.1.3";  I=99.7;              J="2024-09-03" },
    @{ Row=22; A="Intel(R) Wi-Fi 6 AX201 160MHz - 22.0.0.6";   B=1833336; C=2376;  D=656;   E=1899; F=1836368; H="22.0.0.6";   I=99.8;              J="2020-09-16" },
    @{ Row=23; A="Intel(R) Wi-Fi 6 AX201 160MHz - 22.170.0.3"; B=1237156; C=2499;  D=583;   E=1951; F=1240238; H="22.170.0.3"; I=99.8;              J="2022-08-28" },
    @{ Row=24; A="Intel(R) Wi-Fi 6 AX201 160MHz - 21.80.0.4";  B=63994;   C=66;    D=41;    E=68;   F=64101;   H="21.80.0.4";  I=99.8;              J="2020-01-29" },
    @{ Row=25; A="Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9";  B=81417;   C=107;   D=25;    E=158;  F=81549;   H="22.80.0.9";  I=99.8;              J="2021-08-18" },
    @{ Row=26; A="Intel(R) Wi-Fi 6 AX201 160MHz - 22.0.1.1";   B=1400496; C=2029;  D=1075;  E=3210; F=1403600; H="22.0.1.1";   I=99.8;              J="2020-09-28" },
    @{ Row=27; A="Intel(R) Wi-Fi 6 AX201 160MHz - 21.80.2.1";  B=393045;  C=453;   D=216;   E=1251; F=393714;  H="21.80.2.1"; I=99.8;               J="2020-02-24" },
    @{ Row=28; A="Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"; B=509717;  C=522;   D=178;   E=590;  F=510417;  H="23.100.0.4"; I=99.90000000000001; J="2024-11-10" },
    @{ Row=29; A="Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3";  B=11362;   C=7;     D=0;     E=30;   F=11369;   H="21.40.1.3";  I=99.90000000000001; J=$null },
    @{ Row=30; A="Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"; B=75454;   C=2;     D=1;     E=79;   F=75457;   H="21.110.3.2"; I=100;               J="2020-08-05" },
    @{ Row=31; A="Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1";  B=56003;   C=0;     D=15;    E=62;   F=56018;   H="21.60.2.1";  I=100;               J="2019-12-14" },
    @{ Row=32; A="Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6";  B=121310;  C=26;    D=14;    E=154;  F=121350;  H="21.70.0.6";  I=100;               J="2020-01-06" },
    @{ Row=33; A="Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1";  B=35363;   C=11;    D=3;     E=81;   F=35377;   H="22.50.1.1";  I=100;               J="2021-04-27" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    # Column G (lowercased adapter name) is unchanged by this update.
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I

    $jCell = $ws.Cells.Item($row, 10)
    if ($null -eq $r.J) {
        $jCell.ClearContents()
    } else {
        # Force text formatting before assignment so date-shaped strings
        # (e.g. "2023-07-25") are stored as literal text, not serial dates.
        $jCell.NumberFormat = "@"
        $jCell.Value = $r.J
        $jCell.Style = "Normal"
    }
}
